# Final Presentations and ADRs
#
# 1. Resize the "TextBox 7" shape on slide 1 (index 3) - keep its
#    position, change its size.
# 2. Append a new blank slide at the end of the deck.

$p = $ppt.ActivePresentation

# --- 1. Resize shape on slide 1 -------------------------------------------
$slide1 = $p.Slides.Item(1)
$shape = $slide1.Shapes.Item(3)

# EMU -> point conversion (914400 EMU/in, 72 pt/in => 12700 EMU/pt).
# Nudge by half an EMU before dividing so that the host's point -> EMU
# truncation lands back on the exact target EMU value.
$emuPerPt = 12700.0
$shape.Width = (3737113 + 0.5) / $emuPerPt
$shape.Height = (1569660 + 0.5) / $emuPerPt

# --- 2. Add a new blank slide at the end of the deck -----------------------
$newIndex = $p.Slides.Count + 1
$newSlide = $p.Slides.Add($newIndex, 12)
